$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right went from 5 to 4, Wrong went from -1 to -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right went from 95 to 76, Wrong went from -9 to -18
$ws.Range("B12").Value = 76
$ws.Range("C12").Value = -18
$ws.Range("E12").Value = "58 / 112"
